$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds a new "Bill ID" query plus its FC/EC values. In the source
# workbook every data cell is stored as a shared-string (text), even though
# the literal content ("000965482078", "8091", "6741.63") looks numeric
# (note the significant leading zeros on the bill id). Assigning such a
# literal straight to Range.Value would make Excel auto-detect it as a
# number and we'd lose the leading zeros / change the stored cell type.
#
# So each new value is first written with a leading apostrophe into a
# scratch cell (forcing Excel to keep it as literal text), copied, and
# pasted as values-only into the destination - PasteSpecial(values) carries
# over the text content/type but not the quote-prefix formatting, so the
# destination cells keep the workbook's original (default) style. The
# scratch cell is cleared afterwards so it leaves no trace.
$scratch = $ws.Range("Z1")

function Set-TextValue($rangeAddress, $text) {
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $ws.Range($rangeAddress).PasteSpecial("xlPasteValues")
}

Set-TextValue "A2" "000965482078"
Set-TextValue "B2" "8091"
Set-TextValue "C2" "6741.63"
Set-TextValue "D2" "6741.63"

$scratch.Clear()
